{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\n    \"Play Dungeon Quest for Free - Review of Gameplay, Features and Strategies\",\n    \"Play Dungeon Quest Free: Unique Layout and Lucrative Bonus Features\"\n  ],\n  [\n    \"Unique layout and plenty of opportunities to win\",\n    \"Unique layout and plenty of opportunities to discover treasures\"\n  ],\n  [\n    \"Simple to play with an auto-spin function\",\n    \"Simple gameplay with auto-spin function for convenience\"\n  ],\n  [\n    \"Cartoonish design with adventurous music\",\n    \"Pleasing visuals and sound design\"\n  ],\n  [\n    \"Reasonably standard RTP of 96.27%\",\n    \"Standard RTP compared to other slot games\"\n  ],\n  [\n    \"Lack of flashy 3D effects may not appeal to all players\",\n    \"Lack of flashy 3D effects for those seeking a more immersive experience\"\n  ],\n  [\n    \"Discover the treasures of the Dwarf mines in Dungeon Quest, a unique and entertaining game with bonus features. Play now for free and learn tips to win.\",\n    \"Discover treasures in Dungeon Quest free, with unique gameplay and bonus features.\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Dungeon Quest for Free - Review of Gameplay, Features and Strategies\"; Replace = \"Play Dungeon Quest Free: Unique Layout and Lucrative Bonus Features\" },\n    @{ Find = \"Unique layout and plenty of opportunities to win\"; Replace = \"Unique layout and plenty of opportunities to discover treasures\" },\n    @{ Find = \"Simple to play with an auto-spin function\"; Replace = \"Simple gameplay with auto-spin function for convenience\" },\n    @{ Find = \"Cartoonish design with adventurous music\"; Replace = \"Pleasing visuals and sound design\" },\n    @{ Find = \"Reasonably standard RTP of 96.27%\"; Replace = \"Standard RTP compared to other slot games\" },\n    @{ Find = \"Lack of flashy 3D effects may not appeal to all players\"; Replace = \"Lack of flashy 3D effects for those seeking a more immersive experience\" },\n    @{ Find = \"Discover the treasures of the Dwarf mines in Dungeon Quest, a unique and entertaining game with bonus features. Play now for free and learn tips to win.\"; Replace = \"Discover treasures in Dungeon Quest free, with unique gameplay and bonus features.\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
